$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting (bold, borders, centered) from H1, then set new text
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-26: I = 1, J = value from column H
for ($r = 2; $r -le 26; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
